# Generate Report for Handback
# Updates the timestamp strings recorded on the handback-status report
# to reflect the latest run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet, and the matching
# "Correspond Handoff Datetime" on the de-de sheet share the same value.
$wsOverview.Range("G2").Value() = "2016-08-20 01:08:09"
$wsDeDe.Range("H2").Value() = "2016-08-20 01:08:09"

# zh-cn handoff/handback timestamps.
$wsZhCn.Range("H2").Value() = "2016-08-20 01:08:01"
$wsZhCn.Range("K2").Value() = "2016-08-20 01:08:29"

# de-de handback timestamp.
$wsDeDe.Range("K2").Value() = "2016-08-20 01:08:35"
